$d = $word.ActiveDocument

# 1) Remove the "Meta description" paragraph from the top of the document
#    (it sat right after the H1 title).
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2) Insert a new bold paragraph repeating the page title just above the
#    closing (italic) image-prompt paragraph.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$titleParaXml = "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Book of Souls II: El Dorado for Free - Review</w:t></w:r></w:p><w:p/></w:body></w:document>"
$insertPoint.InsertXML($titleParaXml)

# InsertXML leaves a spare empty paragraph behind the new title paragraph
# (before the old closing paragraph) - drop it.
$spareIndex = $d.Paragraphs.Count - 1
$d.Paragraphs($spareIndex).Range.Delete()

# 3) Swap the text of the final (italic) paragraph from the old
#    image-generation prompt to the meta-description copy.
$oldPrompt = 'Create a feature image fitting the game "Book of Souls II: El Dorado". The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior should be holding the Sacred Book of Souls and standing in front of the entrance to the hidden temple. In the background, the temple should be visible with a sense of mystery and adventure.'
$newDescription = "Explore the ancient ruins of El Dorado in Book of Souls II. Enjoy unique features, such as Snake Wilds and two types of free spins for high payout potential. Play now for free."
$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2)
